$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1208.6
$ws.Range("I19").Value = 706.5
$ws.Range("K19").Value = 706.5
$ws.Range("M19").Value = -531.5

$ws.Range("H33").Value = 784.35297
$ws.Range("I33").Value = 222.15384
$ws.Range("K33").Value = 222.15384
$ws.Range("M33").Value = 6.846159999999998

$ws.Range("H113").Value = 9078.362999999999
$ws.Range("I113").Value = 10844
$ws.Range("J113").Value = 5988.5
$ws.Range("K113").Value = 10844
$ws.Range("L113").Value = 5988.5
$ws.Range("M113").Value = -7590
$ws.Range("N113").Value = -12496.5

$ws.Range("H132").Value = 2590.6956
$ws.Range("I132").Value = 2722.2632
$ws.Range("J132").Value = 1965.75
$ws.Range("K132").Value = 8166.7896
$ws.Range("L132").Value = 5897.25
$ws.Range("M132").Value = -5636.7896
$ws.Range("N132").Value = -10957.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 580
$ws.Range("I4").Value = 900
$ws.Range("K4").Value = 900
$ws.Range("M4").Value = -784

$ws.Range("H32").Value = 2336.239
$ws.Range("I32").Value = 2283.4546
$ws.Range("K32").Value = 2283.4546
$ws.Range("M32").Value = -1996.4546

$ws.Range("H45").Value = 10215.19
$ws.Range("I45").Value = 11732.571
$ws.Range("K45").Value = 11732.571
$ws.Range("M45").Value = -11355.571

$ws.Range("H61").Value = 6722.923
$ws.Range("I61").Value = 7376.3335
$ws.Range("K61").Value = 7376.3335
$ws.Range("M61").Value = -7164.3335

$ws.Range("H74").Value = 3477.4119
$ws.Range("I74").Value = 1777.8846
$ws.Range("J74").Value = 9000.875
$ws.Range("K74").Value = 1777.8846
$ws.Range("L74").Value = 9000.875
$ws.Range("M74").Value = -903.8846000000001
$ws.Range("N74").Value = -10748.875

$ws.Range("H77").Value = 3477.4119
$ws.Range("I77").Value = 1777.8846
$ws.Range("J77").Value = 9000.875
$ws.Range("K77").Value = 8889.423000000001
$ws.Range("L77").Value = 45004.375
$ws.Range("M77").Value = -4521.423000000001
$ws.Range("N77").Value = -53740.375

$ws.Range("H102").Value = 5181.381
$ws.Range("I102").Value = 2801.0908
$ws.Range("K102").Value = 2801.0908
$ws.Range("M102").Value = -1179.0908

$ws.Range("H132").Value = 3226.5264
$ws.Range("I132").Value = 1880.9166
$ws.Range("K132").Value = 5642.7498
$ws.Range("M132").Value = -3112.7498

$ws.Range("H136").Value = 6722.923
$ws.Range("I136").Value = 7376.3335
$ws.Range("K136").Value = 22129.0005
$ws.Range("M136").Value = -19579.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3939.6667
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H134").Value = 2351.6428
$ws.Range("I134").Value = 2183.725
$ws.Range("J134").Value = 5710
$ws.Range("K134").Value = 6551.174999999999
$ws.Range("L134").Value = 17130
$ws.Range("M134").Value = -4016.174999999999
$ws.Range("N134").Value = -22200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 844.4358999999999
$ws.Range("I22").Value = 828.4
$ws.Range("K22").Value = 828.4
$ws.Range("M22").Value = -478.4

$ws.Range("H50").Value = 18200

$ws.Range("H51").Value = 37500

$ws.Range("H59").Value = 39166.668
$ws.Range("J59").Value = 39166.668
$ws.Range("L59").Value = 39166.668
$ws.Range("N59").Value = -41456.668

$ws.Range("H61").Value = 37500

$ws.Range("H105").Value = 18420.285
$ws.Range("J105").Value = 4950
$ws.Range("L105").Value = 4950
$ws.Range("N105").Value = -8444

$ws.Range("H107").Value = 8623.1
$ws.Range("I107").Value = 12149.9
$ws.Range("K107").Value = 12149.9
$ws.Range("M107").Value = -10229.9

$ws.Range("H132").Value = 31618.5
$ws.Range("I132").Value = 2081.7778
$ws.Range("J132").Value = 84784.60000000001
$ws.Range("K132").Value = 6245.3334
$ws.Range("L132").Value = 254353.8
$ws.Range("M132").Value = -3715.3334
$ws.Range("N132").Value = -259413.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 591.25
$ws.Range("I2").Value = 400.25
$ws.Range("K2").Value = 2401.5
$ws.Range("M2").Value = -2288.5

$ws.Range("H26").Value = 222.66667
$ws.Range("I26").Value = 154.7
$ws.Range("K26").Value = 464.1
$ws.Range("M26").Value = -176.1

$ws.Range("H34").Value = 1924739.1
$ws.Range("J34").Value = 2199.6
$ws.Range("L34").Value = 6598.799999999999
$ws.Range("N34").Value = -6766.799999999999

$ws.Range("H38").Value = 1396.5834
$ws.Range("J38").Value = 2920.2
$ws.Range("L38").Value = 8760.599999999999
$ws.Range("N38").Value = -9454.599999999999

$ws.Range("H39").Value = 11812.5
$ws.Range("J39").Value = 32500
$ws.Range("L39").Value = 97500
$ws.Range("N39").Value = -98088

$ws.Range("H44").Value = 1470.0714
$ws.Range("J44").Value = 2016.8572
$ws.Range("L44").Value = 6050.571599999999
$ws.Range("N44").Value = -6846.571599999999

$ws.Range("H55").Value = 4192.6665
$ws.Range("J55").Value = 4472.75
$ws.Range("L55").Value = 13418.25
$ws.Range("N55").Value = -13772.25

$ws.Range("H134").Value = 2035.2
$ws.Range("I134").Value = 2035.2
$ws.Range("K134").Value = 6105.6
$ws.Range("M134").Value = -1035.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 10428.286
$ws.Range("I97").Value = 10933.05
$ws.Range("J97").Value = 333
$ws.Range("K97").Value = 10933.05
$ws.Range("L97").Value = 333
$ws.Range("M97").Value = -10437.05
$ws.Range("N97").Value = -1325

$ws.Range("H132").Value = 4535.091
$ws.Range("I132").Value = 2269.5715
$ws.Range("J132").Value = 8499.75
$ws.Range("K132").Value = 6808.7145
$ws.Range("L132").Value = 25499.25
$ws.Range("M132").Value = -4278.7145
$ws.Range("N132").Value = -30559.25

$ws.Range("H133").Value = 54853.332
$ws.Range("J133").Value = 54853.332
$ws.Range("L133").Value = 54853.332
$ws.Range("N133").Value = -64973.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6667.174
$ws.Range("I93").Value = 7216.95
$ws.Range("K93").Value = 7216.95
$ws.Range("M93").Value = -5968.95

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 60000
$ws.Range("J56").Value = 60000
$ws.Range("L56").Value = 60000
$ws.Range("N56").Value = -61428

$ws.Range("H76").Value = 28000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 28000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 28000
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -28630

$ws.Range("H79").Value = 28000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 28000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 28000
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -30184

$ws.Range("H132").Value = 7559.1963
$ws.Range("I132").Value = 7688.143
$ws.Range("K132").Value = 23064.429
$ws.Range("M132").Value = -20534.429

$ws.Range("H136").Value = 2777.7036
$ws.Range("I136").Value = 2000.12
$ws.Range("K136").Value = 6000.36
$ws.Range("M136").Value = -3450.36
